$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks error: marking scheme for right/wrong answers was updated
# (right: 5 -> 4, wrong: -1 -> -2), which cascades into the Total row and
# the Max score summary text.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "78 / 112"
